# Update "想去人数" (interested-people count) figures on the 展览 and
# 全部类型 sheets, and one figure on 演出, per the gh-pages data refresh.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsShow       = $wb.Worksheets.Item("演出")
$wsAll        = $wb.Worksheets.Item("全部类型")

# 展览 (sheet1)
$wsExhibition.Range("F2").Value  = 4630
$wsExhibition.Range("F3").Value  = 2522
$wsExhibition.Range("F10").Value = 186
$wsExhibition.Range("F11").Value = 176
$wsExhibition.Range("F12").Value = 1732
$wsExhibition.Range("F13").Value = 319
$wsExhibition.Range("F14").Value = 3848
$wsExhibition.Range("F15").Value = 35
$wsExhibition.Range("F16").Value = 258

# 演出 (sheet2)
$wsShow.Range("F5").Value = 10

# 全部类型 (sheet4)
$wsAll.Range("F2").Value  = 4630
$wsAll.Range("F3").Value  = 2522
$wsAll.Range("F12").Value = 186
$wsAll.Range("F13").Value = 176
$wsAll.Range("F15").Value = 10
$wsAll.Range("F16").Value = 1732
$wsAll.Range("F17").Value = 319
$wsAll.Range("F18").Value = 3848
$wsAll.Range("F19").Value = 35
$wsAll.Range("F20").Value = 258
